# "Removed S.No from Sample file"
# The sample-format sheet had a leading "S.No" column (A) followed by the
# real header row (regNumber, firstName, middleName, lastName, ...), plus a
# few blank, pre-formatted rows (2-6) below the headers used only to carry
# date-number-format styling. The edit drops the S.No column entirely
# (shifting every other column one slot to the left) and removes the
# now-pointless blank formatted rows, leaving a clean single header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A ("S.No") - everything else shifts left automatically.
$ws.Columns("A").Delete() | Out-Null

# The blank rows below the header only held left-over number-format
# styling (no real values) - remove them so the sheet is just the header.
$ws.Rows("2:6").Delete() | Out-Null

# Park the selection where the saved file shows it.
$ws.Range("E22").Select() | Out-Null
